$wb = $excel.ActiveWorkbook

# "Correspond Handoff Datetime" (col D) and "Correspond Handback DateTime" (col G)
# on row 2 of each language sheet get refreshed with a newly generated handback
# report timestamp for the 653de435-... file pair.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(2, 4).Value = "2016-02-18 09:52:16"
$wsZhCn.Cells.Item(2, 7).Value = "2016-02-18 09:53:06"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(2, 4).Value = "2016-02-18 09:52:31"
$wsDeDe.Cells.Item(2, 7).Value = "2016-02-18 09:53:28"
